$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 5 and row 6 (A:D)
$ws.Range("A5").Value = "25-03-2025"
$ws.Range("B5").Value = "Gujarat Titans vs Punjab Kings"
$ws.Range("C5").Value = "Gujarat Titans"
$ws.Range("D5").Value = "Gujarat Titans"

$ws.Range("A6").Value = "24-03-2025"
$ws.Range("B6").Value = "Delhi Capitals vs Lucknow Super Giants"
$ws.Range("C6").Value = "Delhi Capitals"
$ws.Range("D6").Value = "Delhi Capitals"
